$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1020.3571
$ws.Range("I32").Value = 112.5
$ws.Range("J32").Value = 1383.5
$ws.Range("K32").Value = 112.5
$ws.Range("L32").Value = 1383.5
$ws.Range("M32").Value = 213.5
$ws.Range("N32").Value = -2035.5

$ws.Range("H33").Value = 208.03847
$ws.Range("I33").Value = 160.86363
$ws.Range("K33").Value = 160.86363
$ws.Range("M33").Value = 68.13637

$ws.Range("H100").Value = 2963.75
$ws.Range("I100").Value = 2844.1667
$ws.Range("J100").Value = 3083.3333
$ws.Range("K100").Value = 2844.1667
$ws.Range("L100").Value = 3083.3333
$ws.Range("M100").Value = -2303.1667
$ws.Range("N100").Value = -4165.3333

$ws.Range("H116").Value = 1654.1333
$ws.Range("I116").Value = 1487.5
$ws.Range("K116").Value = 1487.5
$ws.Range("M116").Value = 1954.5

$ws.Range("H132").Value = 964431.4
$ws.Range("I132").Value = 3598.1353
$ws.Range("K132").Value = 10794.4059
$ws.Range("M132").Value = -8264.4059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4230.27
$ws.Range("I32").Value = 3726.7896
$ws.Range("J32").Value = 13796.4
$ws.Range("K32").Value = 3726.7896
$ws.Range("L32").Value = 13796.4
$ws.Range("M32").Value = -3439.7896
$ws.Range("N32").Value = -14370.4

$ws.Range("H61").Value = 23304252
$ws.Range("I61").Value = 32291746
$ws.Range("J61").Value = 86561.836
$ws.Range("K61").Value = 32291746
$ws.Range("L61").Value = 86561.836
$ws.Range("M61").Value = -32291534
$ws.Range("N61").Value = -86985.836

$ws.Range("H74").Value = 4473266.5
$ws.Range("I74").Value = 5308665.5
$ws.Range("K74").Value = 5308665.5
$ws.Range("M74").Value = -5307791.5

$ws.Range("H77").Value = 4473266.5
$ws.Range("I77").Value = 5308665.5
$ws.Range("K77").Value = 26543327.5
$ws.Range("M77").Value = -26538959.5

$ws.Range("H136").Value = 23304252
$ws.Range("I136").Value = 32291746
$ws.Range("J136").Value = 86561.836
$ws.Range("K136").Value = 96875238
$ws.Range("L136").Value = 259685.508
$ws.Range("M136").Value = -96872688
$ws.Range("N136").Value = -264785.508

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 689.4211
$ws.Range("I94").Value = 660.6
$ws.Range("J94").Value = 797.5
$ws.Range("K94").Value = 660.6
$ws.Range("L94").Value = 797.5
$ws.Range("M94").Value = -209.6
$ws.Range("N94").Value = -1699.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38314.707
$ws.Range("I31").Value = 33009.445
$ws.Range("J31").Value = 41497.867
$ws.Range("K31").Value = 33009.445
$ws.Range("L31").Value = 41497.867
$ws.Range("M31").Value = -32714.445
$ws.Range("N31").Value = -42087.867

$ws.Range("H34").Value = 38314.707
$ws.Range("I34").Value = 33009.445
$ws.Range("J34").Value = 41497.867
$ws.Range("K34").Value = 33009.445
$ws.Range("L34").Value = 41497.867
$ws.Range("M34").Value = -32807.445
$ws.Range("N34").Value = -41901.867

$ws.Range("H99").Value = 1114
$ws.Range("I99").Value = 1182.4
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1182.4
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 315.5999999999999
$ws.Range("N99").Value = -3996

$ws.Range("H105").Value = 1023.4545
$ws.Range("I105").Value = 1023.4545
$ws.Range("K105").Value = 1023.4545
$ws.Range("M105").Value = 723.5454999999999

$ws.Range("H122").Value = 2112
$ws.Range("I122").Value = 1912.25
$ws.Range("J122").Value = 2911
$ws.Range("K122").Value = 5736.75
$ws.Range("L122").Value = 8733
$ws.Range("M122").Value = -3286.75
$ws.Range("N122").Value = -13633

$ws.Range("H126").Value = 1114
$ws.Range("I126").Value = 1182.4
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 3547.2
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -1077.2
$ws.Range("N126").Value = -7940

$ws.Range("H133").Value = 46889.85
$ws.Range("J133").Value = 46889.85
$ws.Range("L133").Value = 46889.85
$ws.Range("N133").Value = -51949.85

$ws.Range("H134").Value = 48397.25
$ws.Range("I134").Value = 4333.6665
$ws.Range("J134").Value = 180588
$ws.Range("K134").Value = 13000.9995
$ws.Range("L134").Value = 541764
$ws.Range("M134").Value = -10465.9995
$ws.Range("N134").Value = -546834

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 316
$ws.Range("I2").Value = 375
$ws.Range("J2").Value = 286.5
$ws.Range("K2").Value = 2250
$ws.Range("L2").Value = 1719
$ws.Range("M2").Value = -2137
$ws.Range("N2").Value = -1945

$ws.Range("H122").Value = 936.23334
$ws.Range("I122").Value = 300.5
$ws.Range("K122").Value = 2704.5
$ws.Range("M122").Value = -254.5

$ws.Range("H136").Value = 3181.5293
$ws.Range("I136").Value = 2635.75
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 7907.25
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -2807.25
$ws.Range("N136").Value = -21200.0001

$ws.Range("H137").Value = 1810.3667
$ws.Range("I137").Value = 793.2222
$ws.Range("J137").Value = 3336.0833
$ws.Range("K137").Value = 2379.6666
$ws.Range("L137").Value = 10008.2499
$ws.Range("M137").Value = 2720.3334
$ws.Range("N137").Value = -20208.2499

$ws.Range("H139").Value = 4037.3
$ws.Range("I139").Value = 1624.7742
$ws.Range("J139").Value = 7973.5264
$ws.Range("K139").Value = 4874.3226
$ws.Range("L139").Value = 23920.5792
$ws.Range("M139").Value = 265.6773999999996
$ws.Range("N139").Value = -34200.5792

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2320
$ws.Range("I126").Value = 1866.6666
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 5599.9998
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -3129.9998
$ws.Range("N126").Value = -13940

$ws.Range("H141").Value = 37486.668
$ws.Range("J141").Value = 37486.668
$ws.Range("L141").Value = 37486.668
$ws.Range("N141").Value = -47846.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3030986
$ws.Range("I46").Value = 4329616
$ws.Range("J46").Value = 850
$ws.Range("K46").Value = 4329616
$ws.Range("L46").Value = 850
$ws.Range("M46").Value = -4329428
$ws.Range("N46").Value = -1226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 39926.332
$ws.Range("J69").Value = 39926.332
$ws.Range("L69").Value = 39926.332
$ws.Range("N69").Value = -41424.332

$ws.Range("H72").Value = 39926.332
$ws.Range("J72").Value = 39926.332
$ws.Range("L72").Value = 119778.996
$ws.Range("N72").Value = -127266.996

$ws.Range("H122").Value = 4250
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4250
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12750
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -17650

$ws.Range("H136").Value = 51594.195
$ws.Range("I136").Value = 35992.31
$ws.Range("J136").Value = 89298.75
$ws.Range("K136").Value = 107976.93
$ws.Range("L136").Value = 267896.25
$ws.Range("M136").Value = -105426.93
$ws.Range("N136").Value = -272996.25
